$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, pushing the old rows 32 (VAUPÉS) and 33 (VICHADA) down
$ws.Rows.Item(32).Insert()

# Row 31 (GUAINÍA) code value changes from 95 to 94
$ws.Cells.Item(31, 1).Value = 94

# New row 32: code 95, department GUAVIARE
$ws.Cells.Item(32, 1).Value = 95
$ws.Cells.Item(32, 2).Value = "GUAVIARE"

# Explicitly set row 2 height so it is persisted as a custom row height
$ws.Rows.Item(2).RowHeight = 15
